$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PlainValue($cellRef, $value) {
    $ws.Range($cellRef).Value = $value
}

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "67.550.34"
Set-TextValue "E2" "  -0.27%  "
Set-TextValue "D3" "3.779.11"
Set-TextValue "E3" "  -0.38%  "
Set-TextValue "E4" "  +0.08%  "
Set-TextValue "D5" "597.69"
Set-TextValue "D6" "164.54"
Set-TextValue "E6" "  -1.53%  "
Set-TextValue "E7" "  -0.14%  "
Set-TextValue "D8" "0.515"
Set-TextValue "E8" "  -1.06%  "
Set-TextValue "E9" "  -1.06%  "
Set-TextValue "E10" "  +0.24%  "
Set-TextValue "E11" "  +1.24%  "
Set-TextValue "E12" "  -2.04%  "
Set-TextValue "D13" "35.52"
Set-TextValue "E13" "  -1.49%  "
Set-TextValue "D14" "4.414.10"
Set-TextValue "E14" "  -0.33%  "
Set-TextValue "D15" "3.784.47"
Set-TextValue "E15" "  +0.63%  "
Set-TextValue "D16" "67.591.50"
Set-TextValue "E16" "  -0.16%  "
Set-TextValue "D17" "18.28"
Set-TextValue "E17" "  -0.83%  "
Set-TextValue "E18" "  +1.73%  "
Set-TextValue "D19" "7.01"
Set-TextValue "E19" "  -0.60%  "
Set-TextValue "D20" "460.13"
Set-TextValue "E20" "  +0.12%  "
Set-TextValue "D21" "9.73"
Set-TextValue "E21" "  -2.44%  "
Set-TextValue "D22" "0.695"
Set-TextValue "E22" "  -0.20%  "
Set-TextValue "E23" "  -5.80%  "
Set-TextValue "D24" "82.39"
Set-TextValue "E24" "  -1.21%  "
Set-TextValue "D25" "11.97"
Set-TextValue "E25" "  -0.91%  "
Set-TextValue "D26" "2.09"
Set-TextValue "E26" "  -1.21%  "
Set-TextValue "D27" "0.999"
Set-TextValue "E27" "  -0.10%  "
Set-TextValue "E28" "  -0.42%  "
Set-TextValue "D29" "3.927.45"
Set-TextValue "E29" "  -0.36%  "
Set-TextValue "D30" "7.44"
Set-TextValue "E30" "  +2.96%  "
Set-TextValue "E31" "  -4.54%  "
Set-TextValue "D32" "2.19"
Set-TextValue "E32" "  -2.62%  "
Set-TextValue "D33" "28.98"
Set-TextValue "E33" "  -2.16%  "
Set-TextValue "D34" "0.998"
Set-TextValue "E34" "  -0.17%  "
Set-TextValue "D35" "8.97"
Set-TextValue "E35" "  -1.12%  "
Set-TextValue "E36" "  -1.69%  "
Set-TextValue "E37" "  +0.20%  "
Set-TextValue "D38" "0.986"
Set-TextValue "E38" "  -0.48%  "
Set-TextValue "D39" "3.22"
Set-TextValue "E39" "  -4.77%  "
Set-TextValue "E40" "  -0.81%  "
Set-TextValue "D41" "0.999"
Set-TextValue "E41" "  +0.02%  "
Set-PlainValue "B43" "OKB"
Set-PlainValue "C43" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D43" "47.48"
Set-TextValue "E43" "  -1.18%  "
Set-PlainValue "B44" "Arweave"
Set-PlainValue "C44" "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue "D44" "43.43"
Set-TextValue "E44" "  -1.14%  "
Set-TextValue "E45" "  -0.68%  "
Set-TextValue "E46" "  +1.09%  "
Set-TextValue "E47" "  +0.31%  "
Set-TextValue "E48" "  +8.87%  "
Set-TextValue "D49" "27.27"
Set-TextValue "E49" "  +1.84%  "
Set-PlainValue "B50" "Stacks"
Set-PlainValue "C50" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D50" "1.85"
Set-TextValue "E50" "  +1.47%  "
Set-PlainValue "B51" "Bittensor"
Set-PlainValue "C51" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D51" "392.70"
Set-TextValue "E51" "  +0.52%  "
